$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.002.54'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '1.888.46'
$ws.Range('E3').Value = '  -3.73%  '
$ws.Range('E4').Value = '  -1.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4584'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3945'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.30'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08208'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').Value = '1.919.45'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.331'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.988'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001058'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06566'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.653'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('D23').Value = '28.010.92'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.310'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = '2.139.47'
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.113'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.685'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '124.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09541'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9591'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.474'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.629'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.471'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('E37').Value = '  -2.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.248'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06110'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.631'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1892'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.314'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5813'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('E48').Value = '  -3.59%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06886'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '110.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.52%  '
